$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 3: extend the thick-bottom-border blank row (style copied from M3) ---
$ws.Range("M3").Copy()
$ws.Range("N3:O3").PasteSpecial(-4122)

# --- Row 4: extend the year header row (style copied from M4) ---
$ws.Range("M4").Copy()
$ws.Range("N4:O4").PasteSpecial(-4122)
$ws.Range("N4").Value = 2021
$ws.Range("O4").Value = 2022

# --- Row 5 (national total, bold + 0.0 number format) ---
$ws.Range("N5").Value = 40.007977647471066
$ws.Range("N5").Font.Bold = $true
$ws.Range("N5").NumberFormat = "0.0"
$ws.Range("N5").Copy()
$ws.Range("O5").PasteSpecial(-4122)
$ws.Range("O5").Value = 42.620582506455563

# --- Rows 6-13 (regular data rows, 0.0 number format) ---
$ws.Range("N6").Value = 5.7072514621689896
$ws.Range("N6").NumberFormat = "0.0"
$ws.Range("N6").Copy()
$ws.Range("N7:N13").PasteSpecial(-4122)
$ws.Range("O6:O13").PasteSpecial(-4122)

$ws.Range("O6").Value = 8.1443914479075037
$ws.Range("N7").Value = 8.9893229854028949
$ws.Range("O7").Value = 10.715961386284755
$ws.Range("N8").Value = 66.307512472824584
$ws.Range("O8").Value = 81.977461999426666
$ws.Range("N9").Value = 23.475213049310256
$ws.Range("O9").Value = 29.828871240443185
$ws.Range("N10").Value = 9.8045372040896162
$ws.Range("O10").Value = 9.7218425128664112
$ws.Range("N11").Value = 9.3737779268960448
$ws.Range("O11").Value = 8.6167819403064012
$ws.Range("N12").Value = 70.457032471318783
$ws.Range("O12").Value = 69.915337594090886
$ws.Range("N13").Value = 98.411252120183207
$ws.Range("O13").Value = 99.08571752721997

# --- Row 14 (bottom border + 0.0 number format) ---
$ws.Range("K3").Copy()
$ws.Range("N14").PasteSpecial(-4122)
$ws.Range("N14").Value = 63.900563564170795
$ws.Range("N14").NumberFormat = "0.0"
$ws.Range("N14").Copy()
$ws.Range("O14").PasteSpecial(-4122)
$ws.Range("O14").Value = 64.805252627098838

# --- Selection as recorded at save time ---
$ws.Range("P8").Select()
